# Auto-generated Excel COM-interop edit script
# Applies numeric value updates to the Bahamut_Profits leve-profit tables
# across sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 123
$ws.Range("H123").Value = 20900.723
$ws.Range("J123").Value = 20900.723
$ws.Range("L123").Value = 20900.723
$ws.Range("N123").Value = -30700.723

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 7363.1406
$ws.Range("I32").Value = 3427.131
$ws.Range("J32").Value = 31372.8
$ws.Range("K32").Value = 3427.131
$ws.Range("L32").Value = 31372.8
$ws.Range("M32").Value = -3140.131
$ws.Range("N32").Value = -31946.8

# Row 102
$ws.Range("H102").Value = 4728.625
$ws.Range("I102").Value = 4728.625
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 4728.625
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -3106.625
$ws.Range("N102").ClearContents()

# Row 132
$ws.Range("H132").Value = 1254.7021
$ws.Range("I132").Value = 800.3333
$ws.Range("J132").Value = 2325.7144
$ws.Range("K132").Value = 2400.9999
$ws.Range("L132").Value = 6977.1432
$ws.Range("M132").Value = 129.0001000000002
$ws.Range("N132").Value = -12037.1432

# Row 133
$ws.Range("H133").Value = 34550
$ws.Range("J133").Value = 34550
$ws.Range("L133").Value = 34550
$ws.Range("N133").Value = -39610

$ws = $wb.Worksheets.Item("BSM")
# Row 99
$ws.Range("H99").Value = 2304.6191
$ws.Range("I99").Value = 2211.0557
$ws.Range("J99").Value = 2866
$ws.Range("K99").Value = 2211.0557
$ws.Range("L99").Value = 2866
$ws.Range("M99").Value = -713.0556999999999
$ws.Range("N99").Value = -5862

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2572.0908
$ws.Range("I31").Value = 2544.484
$ws.Range("K31").Value = 2544.484
$ws.Range("M31").Value = -2249.484

# Row 34
$ws.Range("H34").Value = 2572.0908
$ws.Range("I34").Value = 2544.484
$ws.Range("K34").Value = 2544.484
$ws.Range("M34").Value = -2342.484

# Row 38
$ws.Range("H38").Value = 3300
$ws.Range("I38").Value = 3300
$ws.Range("K38").Value = 3300
$ws.Range("M38").Value = -2923

# Row 46
$ws.Range("H46").Value = 3300
$ws.Range("I46").Value = 3300
$ws.Range("K46").Value = 3300
$ws.Range("M46").Value = -3089

$ws = $wb.Worksheets.Item("CUL")
# Row 49
$ws.Range("H49").Value = 850.6
$ws.Range("I49").Value = 850.6
$ws.Range("K49").Value = 2551.8
$ws.Range("M49").Value = -2395.8

# Row 63
$ws.Range("H63").Value = 4736.25
$ws.Range("I63").Value = 3945
$ws.Range("J63").Value = 5000
$ws.Range("K63").Value = 11835
$ws.Range("L63").Value = 15000
$ws.Range("M63").Value = -11086
$ws.Range("N63").Value = -16498

# Row 66
$ws.Range("H66").Value = 4736.25
$ws.Range("I66").Value = 3945
$ws.Range("J66").Value = 5000
$ws.Range("K66").Value = 35505
$ws.Range("L66").Value = 45000
$ws.Range("M66").Value = -31761
$ws.Range("N66").Value = -52488

# Row 76
$ws.Range("H76").Value = 4966.6665
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 4966.6665
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 14899.9995
$ws.Range("M76").ClearContents()
$ws.Range("N76").Value = -15665.9995

# Row 79
$ws.Range("H79").Value = 4966.6665
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 4966.6665
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 14899.9995
$ws.Range("M79").ClearContents()
$ws.Range("N79").Value = -17551.9995

# Row 87
$ws.Range("H87").Value = 14233.333
$ws.Range("I87").Value = 3940
$ws.Range("K87").Value = 11820
$ws.Range("M87").Value = -10572

# Row 90
$ws.Range("H90").Value = 14233.333
$ws.Range("I90").Value = 3940
$ws.Range("K90").Value = 35460
$ws.Range("M90").Value = -29220

# Row 113
$ws.Range("H113").Value = 759.41174
$ws.Range("I113").Value = 560
$ws.Range("J113").Value = 786
$ws.Range("K113").Value = 1680
$ws.Range("L113").Value = 2358
$ws.Range("M113").Value = 490
$ws.Range("N113").Value = -6698

# Row 131
$ws.Range("H131").Value = 780.6667
$ws.Range("J131").Value = 809.087
$ws.Range("L131").Value = 2427.261
$ws.Range("N131").Value = -12507.261

# Row 132
$ws.Range("H132").Value = 616.36365
$ws.Range("I132").Value = 550
$ws.Range("J132").Value = 1280
$ws.Range("K132").Value = 4950
$ws.Range("L132").Value = 11520
$ws.Range("M132").Value = -2420
$ws.Range("N132").Value = -16580

# Row 137
$ws.Range("H137").Value = 60475.168
$ws.Range("I137").Value = 2562
$ws.Range("J137").Value = 132866.62
$ws.Range("K137").Value = 7686
$ws.Range("L137").Value = 398599.86
$ws.Range("M137").Value = -2586
$ws.Range("N137").Value = -408799.86

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 6333.3335
$ws.Range("I70").Value = 4500
$ws.Range("J70").Value = 10000
$ws.Range("K70").Value = 4500
$ws.Range("L70").Value = 10000
$ws.Range("M70").Value = -4230
$ws.Range("N70").Value = -10540

# Row 73
$ws.Range("H73").Value = 6333.3335
$ws.Range("I73").Value = 4500
$ws.Range("J73").Value = 10000
$ws.Range("K73").Value = 4500
$ws.Range("L73").Value = 10000
$ws.Range("M73").Value = -3564
$ws.Range("N73").Value = -11872

# Row 107
$ws.Range("H107").Value = 293.53333
$ws.Range("I107").Value = 411.55554
$ws.Range("J107").Value = 116.5
$ws.Range("K107").Value = 411.55554
$ws.Range("L107").Value = 116.5
$ws.Range("M107").Value = 1508.44446
$ws.Range("N107").Value = -3956.5

# Row 132
$ws.Range("H132").Value = 2775.389
$ws.Range("I132").Value = 2613.739
$ws.Range("J132").Value = 3061.3845
$ws.Range("K132").Value = 7841.217000000001
$ws.Range("L132").Value = 9184.1535
$ws.Range("M132").Value = -5311.217000000001
$ws.Range("N132").Value = -14244.1535

$ws = $wb.Worksheets.Item("LTW")
# Row 133
$ws.Range("H133").Value = 42625
$ws.Range("J133").Value = 42625
$ws.Range("L133").Value = 42625
$ws.Range("N133").Value = -47685

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 1702.2122
$ws.Range("I132").Value = 1098.421
$ws.Range("K132").Value = 3295.263
$ws.Range("M132").Value = -765.2629999999999

# Row 136
$ws.Range("H136").Value = 1415.3611
$ws.Range("I136").Value = 1334.9395
$ws.Range("J136").Value = 2300
$ws.Range("K136").Value = 4004.8185
$ws.Range("L136").Value = 6900
$ws.Range("M136").Value = -1454.8185
$ws.Range("N136").Value = -12000

Write-Host "Applied Bahamut_Profits updates"